# The edit swaps the contents of the data columns between row 23 and row
# 24 of the active sheet (two species records traded places). Columns
# C, D, P, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY hold the
# same value in both rows already, so they are intentionally left
# untouched; K and N are likewise identical (blank) in both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","E","F","G","H","I","J","L","M","Q","R","S","AF")

# Capture the "before" values of both rows first, since we are about to
# overwrite them.
$row23 = @{}
$row24 = @{}
foreach ($col in $cols) {
    $row23[$col] = $ws.Range($col + "23").Value2
    $row24[$col] = $ws.Range($col + "24").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "23").Value2 = $row24[$col]
    $ws.Range($col + "24").Value2 = $row23[$col]
}

# Column I holds a value that looks numeric ("1") but must stay text;
# force a text number format on the cell that ends up holding it (matches
# typing an apostrophe-prefixed value in the Excel UI).
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value2 = $row23["I"]

# A few of the swapped-in cells are blank but were present (empty) cells
# in their source row rather than fully absent ones. Touching the number
# format (without changing the already-blank value) keeps them present in
# the saved sheet instead of letting the blank write above drop them.
$ws.Range("I23").NumberFormat = "General"
$ws.Range("J23").NumberFormat = "General"
$ws.Range("L24").NumberFormat = "General"
